# Minor updates in the case cities data form.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Move the "Our area is ..." notes from the population-data Notes column (H)
# to the travel-data Notes column (W) for each populated data row (3-9),
# clearing the old H value and replacing whatever was previously in W.

$ws.Cells.Item(3, 23).Value = "Our area is Greater Accra Metropolitan Area, formed by 12 districts: Accra Metropolitan District, Tema Metropolis, Adenta, Ga East, Ga West, Ga South, Ga Central, La Nkwantang-Madina, Ledzokuku-Krowor, Ashaiman, Kpone-Katamanso, La Dade-Kotopon. Dataset must be filter in for urban Greater Accra region."
$ws.Cells.Item(3, 8).Value = ""

$ws.Cells.Item(4, 23).Value = "Our area is Sao Paulo Metropolitan Area, formed by 39 cities: Caieiras, Cajamar, Francisco Morato, Franco da Rocha, Mairiporã, Arujá, Biritiba-Mirim, Ferraz de Vasconcelos, Guararema, Guarulhos, Itaquaquecetuba, Mogi das Cruzes, Poá, Salesópolis, Santa Isabel, Suzano, Diadema, Mauá, Ribeirão Pires, Rio Grande da Serra, Santo André, São Bernardo do Campo, São Caetano do Sul, Cotia, Embu das Artes, Embu-Guaçu, Itapecerica da Serra, Juquitiba, São Lourenço da Serra, Taboão da Serra, Vargem Grande Paulista, Barueri, Carapicuíba, Itapevi, Jandira, Osasco, Pirapora do Bom Jesus, Santana de Parnaíba, São Paulo."
$ws.Cells.Item(4, 8).Value = ""

$ws.Cells.Item(5, 23).Value = "Our area is Sao Paulo Metropolitan Area, formed by 34 cities: Baldim, Belo Horizonte, Betim, Brumadinho, Caeté, Capim Branco, Confins, Contagem, Esmeraldas, Florestal, Ibirité, Igarapé, Itaguara, Itatiaiuçu, Jaboticatubas, Nova União, Juatuba, Lagoa Santa, Mário Campos, Mateus Leme, Matozinhos, Nova Lima, Pedro Leopoldo, Raposos, Ribeirão das Neves, Rio Acima, Rio Manso, Sabará, Santa Luzia, São Joaquim de Bicas, São José da Lapa, Sarzedo, Taquaraçu de Minas, Vespasiano"
$ws.Cells.Item(5, 8).Value = ""

$ws.Cells.Item(6, 23).Value = "Our area is Bogota D.C. only. Last census done in 2005"
$ws.Cells.Item(6, 8).Value = ""

$ws.Cells.Item(7, 23).Value = "Our area is Valley of Mexico Metropolitan Area, formed by 61 cities: Acolman, Amecameca, Apaxco, Atenco, Atizapán de Zaragoza, Atlautla, Axapusco, Ayapango, Chalco, Chiautla, Chicoloapan, Chiconcuac, Chimalhuacán, Coacalco de Berriozábal, Cocotitlán, Coyotepec, Cuautitlán, Cuautitlán Izcalli, Ecatepec de Morelos, Ecatzingo, Huehuetoca, Hueypoxtla, Huixquilucan, Isidro Fabela, Ixtapaluca, Jaltenco, Jilotzingo, Juchitepec, La Paz, Melchor Ocampo, Mexico City, Naucalpan de Juárez, Nextlalpan, Nezahualcóyotl, Nicolás Romero, Nopaltepec, Otumba, Ozumba, Papalotla, San Martín de las Pirámides, Tecámac, Temamatla, Temascalapa, Tenango del Aire, Teoloyucan, Teotihuacán, Tepetlaoxtoc, Tepetlixpa, Tepotzotlán, Tequixquiac, Texcoco, Tezoyuca, Tizayuca, Tlalmanalco, Tlalnepantla de Baz, Tonanitla, Tultepec, Tultitlán, Valle de Chalco Solidaridad, Villa del Carbón, Zumpango."
$ws.Cells.Item(7, 8).Value = ""

$ws.Cells.Item(8, 23).Value = "Our area is Gran Buenos Aires, formed by 24 cities: Almirante Brown, Avellaneda, Berazategui, Esteban Echeverría, Ezeiza, Florencio Varela, General San Martín, Hurlingham, Ituzaingó, José C. Paz, La Matanza, Lanús, Lomas de Zamora, Malvinas Argentinas, Merlo, Moreno, Morón, Quilmes, San Fernando, San Isidro, San Miguel, Tigre, Tres de Febrero, Vicente López. Sex- and age-specific population projections for the Buenos Aires province."
$ws.Cells.Item(8, 8).Value = ""

$ws.Cells.Item(9, 23).Value = "Our area is Santiago Metropolitan Area, formed by 52 cities:  Alhué, Buin, Calera de Tango, Cerrillos, Cerro Navia, Colina, Conchalí, Curacaví, El Bosque, El Monte, Estación Central, Huechuraba, Independencia, Isla de Maipo, La Cisterna, La Florida, La Granja, La Pintana, La Reina, Lampa, Las Condes, Lo Barnechea, Lo Espejo, Lo Prado, Macul, Maipú, María Pinto, Melipilla, Ñuñoa, Padre Hurtado, Paine, Pedro Aguirre Cerda, Peñaflor, Peñalolén, Pirque, Providencia, Pudahuel, Puente Alto, Quilicura, Quinta Normal, Recoleta, Renca, San Bernardo, San Joaquín, San José de Maipo, San Miguel, San Pedro, San Ramón, Santiago, Talagante, Til Til, Vitacura"
$ws.Cells.Item(9, 8).Value = ""

# Update the view state: zoom level, frozen pane top-left cell, and the active selection.
$sheetView = $ws.Application.ActiveWindow
$sheetView.Zoom = 70
$ws.Range("C3").Select()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("C3").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1:B1").Select()
